$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the sheet (tab) name to reflect the new "through" date
$ws.Name = "Through 2022-11-03"

# Update the row label for November to reflect the new "through" date
$ws.Range("A12").Value = "November (through 11-03)"

# Correct a value in the October row (2022 column)
$ws.Range("I11").Value = 124

# Update the November row (row 12) with new data for years 2016-2022 (columns C-I)
$ws.Range("C12").Value = 6
$ws.Range("D12").Value = 10
$ws.Range("E12").Value = 14
$ws.Range("F12").Value = 4
$ws.Range("G12").Value = 21
$ws.Range("H12").Value = 18
$ws.Range("I12").Value = 9

# Update the Total row (row 13) with new totals for years 2016-2022 (columns C-I)
$ws.Range("C13").Value = 492
$ws.Range("D13").Value = 720
$ws.Range("E13").Value = 629
$ws.Range("F13").Value = 486
$ws.Range("G13").Value = 1078
$ws.Range("H13").Value = 1459
$ws.Range("I13").Value = 1409
